$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Insert a new row before row 18 (the existing data rows 18-21 and the
# spacer row 22 all shift down by one; formulas auto-adjust).
$ws.Rows.Item(18).Insert()

# The freshly inserted row doesn't carry the bordered table formatting, so
# clone it from the row directly below (which still has the original
# look of the data rows). Column A (plain centered/bordered look) is a
# closer match for every column of the new row than the money columns.
$ws.Range("A19").Copy()
$ws.Range("A18:E18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- New row 18: the "PDU beda warna" line item ---
$ws.Range("A18").Value = 1
$ws.Range("B18").Value = "24 MARET 2021"
$ws.Range("C18").Value = "110 SET"
$ws.Range("D18").Value = "-"
$ws.Range("E18").Value = "-"

$ws.Range("B18:C18").HorizontalAlignment = -4131

# Row 19 (previously row 18) keeps its old value/date but now needs a
# left-aligned variant of its style.
$ws.Range("B19").HorizontalAlignment = -4131

# Renumber the "NO" column for the rows that followed.
$ws.Range("A19").Value = 2
$ws.Range("A20").Value = 3
$ws.Range("A21").Value = 4

Write-Host "done"
